$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '59.739.15'
$ws.Range("E2").Value = '  +2.78%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.419.82'
$ws.Range("E3").Value = '  +3.24%  '

$ws.Range("E4").Value = '  -0.12%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '554.44'
$ws.Range("E5").Value = '  +2.58%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '137.75'
$ws.Range("E6").Value = '  +2.06%  '

$ws.Range("E7").Value = '  -0.16%  '

$ws.Range("E8").Value = '  +1.40%  '

$ws.Range("E9").Value = '  +5.12%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '5.80'
$ws.Range("E10").Value = '  +3.27%  '

$ws.Range("E11").Value = '  +1.77%  '

$ws.Range("E12").Value = '  -1.98%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '24.68'
$ws.Range("E13").Value = '  +3.63%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.851.25'
$ws.Range("E14").Value = '  +3.23%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '59.631.55'
$ws.Range("E15").Value = '  +2.62%  '

$ws.Range("E16").Value = '  +4.78%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.405.15'
$ws.Range("E17").Value = '  +2.99%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.35'
$ws.Range("E18").Value = '  +6.22%  '

$ws.Range("E19").Value = '  +4.73%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '334.67'
$ws.Range("E20").Value = '  +0.65%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.95'
$ws.Range("E21").Value = '  +3.49%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.999'
$ws.Range("E22").Value = '  -0.05%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '64.62'
$ws.Range("E23").Value = '  +3.02%  '

$ws.Range("E24").Value = '  +1.07%  '

$ws.Range("E25").Value = '  +2.31%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  -0.11%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.38'
$ws.Range("E27").Value = '  -1.13%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0₃0792'
$ws.Range("E28").Value = '  +7.81%  '

$ws.Range("E29").Value = '  +3.25%  '

$ws.Range("B30").Value = 'Aptos'
$ws.Range("C30").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.29'
$ws.Range("E30").Value = '  +3.15%  '

$ws.Range("B31").Value = 'Monero'
$ws.Range("C31").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '170.20'
$ws.Range("E31").Value = '  +0.22%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '18.72'
$ws.Range("E32").Value = '  +1.95%  '

$ws.Range("E33").Value = '  +1.20%  '

$ws.Range("E34").Value = '  -0.02%  '

$ws.Range("E35").Value = '  +5.48%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.26'
$ws.Range("E36").Value = '  +0.43%  '

$ws.Range("E37").Value = '  +0.15%  '

$ws.Range("E38").Value = '  -0.76%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '40.12'
$ws.Range("E39").Value = '  +2.53%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.419'
$ws.Range("E40").Value = '  +11.04%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '313.72'
$ws.Range("E41").Value = '  +9.19%  '

$ws.Range("E42").Value = '  +3.16%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '142.44'

$ws.Range("E44").Value = '  +2.69%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0523'
$ws.Range("E45").Value = '  +4.30%  '

$ws.Range("B46").Value = 'InjectiveProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '19.22'
$ws.Range("E46").Value = '  +0.56%  '

$ws.Range("B47").Value = 'Mantle'
$ws.Range("C47").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.571'
$ws.Range("E47").Value = '  +1.75%  '

$ws.Range("B48").Value = 'Polygon'
$ws.Range("C48").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.404'
$ws.Range("E48").Value = '  +4.98%  '

$ws.Range("E49").Value = '  +3.12%  '

$ws.Range("E50").Value = '  -0.25%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.62'
$ws.Range("E51").Value = '  +5.64%  '
